$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Shared-string / cell value edits
# ---------------------------------------------------------------------------

# "Cont adminstrator" sheet: school name + email swapped to "Moira"
$wsCont = $wb.Worksheets.Item("Cont adminstrator")
$wsCont.Range("B15").Value = "Scoala particulara Moira"
$wsCont.Range("C15").Value = "moirascoala@automation.33mail.com"

# "Receptie" sheet: three receptionist emails tweaked
$wsRec = $wb.Worksheets.Item("Receptie")
$wsRec.Range("B2").Value = "steluta2a1@staffcalendis.33mail.com"
$wsRec.Range("B3").Value = "steluta3131a@staffcalendis.33mail.com"
$wsRec.Range("B4").Value = "steluta441a@staffcalendis.33mail.com"

# "Angajati" sheet: four employee emails tweaked
$wsAng = $wb.Worksheets.Item("Angajati")
$wsAng.Range("B2").Value = "elenas1114@staffcalendis.33mail.com"
$wsAng.Range("B3").Value = "komornic1124@staffcalendis.33mail.com"
$wsAng.Range("B4").Value = "ovidius134@staffcalendis.33mail.com"
$wsAng.Range("B5").Value = "sdroses4@staffcalendis.33mail.com"

# ---------------------------------------------------------------------------
# 2. Row height tweak on "Angajati" row 5 (14.9 -> 13.8, the sheet default)
# ---------------------------------------------------------------------------
$wsAng.Rows.Item(5).RowHeight = 13.8

# ---------------------------------------------------------------------------
# 3. Remove the mailto hyperlink on Angajati!B5 (also renumbers the
#    legacyDrawing relationship id down from rId3 to rId2 as a side effect)
# ---------------------------------------------------------------------------
$wsAng.Range("B5").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 4. Active tab / selection bookkeeping
#    - workbook activeTab moves from "Cont adminstrator" (0) to "Receptie" (1)
#    - "Cont adminstrator" selection stays C15 but loses tabSelected
#    - "Receptie" becomes the selected tab with active cell B6
#    - "Angajati" active cell moves from B13 to B5
# ---------------------------------------------------------------------------
$wsRec.Activate()
$wsRec.Range("B6").Select()

$wsAng.Range("B5").Select()

# re-activate Receptie last so it is the workbook's active/selected sheet
$wsRec.Activate()

# ---------------------------------------------------------------------------
# 5. New auto-generated "_xlnm._FilterDatabase" defined names
#    Each of the three filtered sheets ("Cont adminstrator", "Domenii",
#    "Domenii existente") gains one more name in the "_0_0_..." chain
#    (83 "_0" suffixes instead of the previous longest chain of 82).
# ---------------------------------------------------------------------------
$suffixParts = @()
for ($i = 0; $i -lt 83; $i++) { $suffixParts += "0" }
$suffix = [string]::Join("_", $suffixParts)
$newFilterDbName = "_xlnm._FilterDatabase_" + $suffix

$wsCont.Names.Add($newFilterDbName, "='Cont adminstrator'!`$A`$1:`$A`$19")

$wsDom = $wb.Worksheets.Item("Domenii")
$wsDom.Names.Add($newFilterDbName, "=Domenii!`$A`$4:`$A`$7")

$wsDomEx = $wb.Worksheets.Item("Domenii existente")
$wsDomEx.Names.Add($newFilterDbName, "='Domenii existente'!`$A`$1:`$Q`$15")
